$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing column H header (Num of Trees in millions -> Num of Real Trees in millions)
$ws.Range("H1").Value = "Num of Real Trees in millions"

# New column headers
$ws.Range("I1").Value = "Num of Fake trees in millions"
$ws.Range("J1").Value = "Avg size of home (sq ft)"

# New data for column I (Num of Fake trees in millions), rows 2-17
$fakeTrees = @(9, 9.3, 9.3, 17.4, 11.7, 11.7, 8.2, 9.5, 10.9, 14.7, 13.9, 12.5, 18.6, 21.1, 23.6, 24.4)

# New data for column J (Avg size of home (sq ft)), rows 2-17
$homeSize = @(2349, 2434, 2469, 2521, 2519, 2438, 2392, 2480, 2505, 2598, 2657, 2687, 2640, 2631, 2623, 2322)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $fakeTrees[$i]
    $ws.Cells.Item($row, 10).Value = $homeSize[$i]
}

# Leave the selection where the user's entry would naturally have ended up
[void]$ws.Range("I19").Select()
